$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

$ws.Range("A11").Value = 8
$ws.Range("C11").Value = "positions of large group"
$ws.Range("D11").Value = "Easy"
$ws.Range("F11").Value = "Easy"
$ws.Range("H11").Value = 30

$ws.Range("H19").Select()
